$d = $word.ActiveDocument

# 1) First paragraph ("Лабораторная работа №13. Todolist. ...") gains a
#    left indent of 707 twips (= 707/20 pt) in addition to the existing
#    firstLine indent of 709 twips.
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 707 / 20.0

# 2) The run that hosts the inline picture ("Пример:" figure) becomes
#    "no proofing" (this is what Word stamps on runs/ranges that hold
#    only a drawing), which serializes as <w:rPr><w:noProof/></w:rPr>
#    on that run.
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = $true
